# Update column F ("dSF") values on Sheet1 to reflect repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = -5
$ws.Range("F3").Value = 13
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = -2
$ws.Range("F8").Value = -4
$ws.Range("F9").Value = -1
$ws.Range("F10").Value = 2
$ws.Range("F11").Value = 4
